$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = -0.1654569698432457
$ws.Range("J2").Value = 0.1369481133787689
$ws.Range("K2").Value = -0.4217527083110341
$ws.Range("L2").Value = 2.034547933685609

$ws.Range("I12").Value = 0.1054748360185471
$ws.Range("J12").Value = 0.04603558480375656
$ws.Range("K12").Value = -0.9714502200944235
$ws.Range("L12").Value = 2.200712507002736

$ws.Range("I13").Value = -0.02713269210245115
$ws.Range("J13").Value = 0.1071674939004506
$ws.Range("K13").Value = -0.598798363972762
$ws.Range("L13").Value = 1.708747698783124

$ws.Range("I14").Value = -0.06622842042182196
$ws.Range("J14").Value = 0.09694862478115245
$ws.Range("K14").Value = -0.2725373584899253
$ws.Range("L14").Value = 1.561925405595501

$ws.Range("I15").Value = -0.06970951543867017
$ws.Range("J15").Value = 0.08317152937608577
$ws.Range("K15").Value = -0.2578607092615042
$ws.Range("L15").Value = 1.805744642355887

$ws.Range("I17").Value = -0.3616129562468736
$ws.Range("J17").Value = 0.1421304991736817
$ws.Range("K17").Value = 0.2505153130106562
$ws.Range("L17").Value = 1.960801187805391
